$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("explanations")
$ws2.Range("A8:G8").Copy()
$ws2.Range("A33:G33").PasteSpecial(-4122)
Write-Output "done"
